# The deck originally has its "Integral" (Red Violet) design applied to the
# main slide master/theme (theme1.xml) while a plain built-in "Office Theme"
# colour set sits unused on the notes master (theme2.xml).
#
# The author simply switched the presentation's applied Design over to the
# built-in "Office Theme" colours (Design tab -> Themes/Variants gallery).
# That swaps which palette is "live" on the slide master: the 12 theme
# colours that used to read Office-Theme blue/orange/grey now become the
# live Red-Violet/Integral colours, and vice versa.
#
# Reproduce that by rewriting the 12 theme colours on the live theme
# (reachable via ThemeColorScheme on the slide master/slides) to the
# standard Office Theme palette.

function RGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Index order exposed by ThemeColorScheme:
#  1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
#  9 accent5, 10 accent6, 11 hlink, 12 folHlink
$tcs.Item(1).RGB  = RGB 0x00 0x00 0x00   # dk1
$tcs.Item(2).RGB  = RGB 0xFF 0xFF 0xFF   # lt1
$tcs.Item(3).RGB  = RGB 0x44 0x54 0x6A   # dk2
$tcs.Item(4).RGB  = RGB 0xE7 0xE6 0xE6   # lt2
$tcs.Item(5).RGB  = RGB 0x5B 0x9B 0xD5   # accent1
$tcs.Item(6).RGB  = RGB 0xED 0x7D 0x31   # accent2
$tcs.Item(7).RGB  = RGB 0xA5 0xA5 0xA5   # accent3
$tcs.Item(8).RGB  = RGB 0xFF 0xC0 0x00   # accent4
$tcs.Item(9).RGB  = RGB 0x44 0x72 0xC4   # accent5
$tcs.Item(10).RGB = RGB 0x70 0xAD 0x47   # accent6
$tcs.Item(11).RGB = RGB 0x05 0x63 0xC1   # hlink
$tcs.Item(12).RGB = RGB 0x95 0x4F 0x72   # folHlink
